# Gigantic update to prepare for test server.
# Applies changes to the "StatDef" worksheet (xl/worksheets/sheet1.xml):
#  - tweak stats on a couple of existing monster rows (289, 290)
#  - switch MonsterAiType from AiAggressive to AiAngry on two rows (295, 299)
#  - append six new boss-monster rows (320-325: Succubus, Incubus, Observation,
#    Retribution, Solace, Shelter)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatDef")

# --- Row 289 (existing monster) stat tweaks ---
$ws.Cells.Item(289,5).Value = 100     # E289 HP 150 -> 100
$ws.Cells.Item(289,9).Value = 70      # I289 Dex 100 -> 70
$ws.Cells.Item(289,10).Value = 30     # J289 Agi 100 -> 30
$ws.Cells.Item(289,11).Value = 100    # K289 Luk 170 -> 100
$ws.Cells.Item(289,12).Value = 40     # L289 Attack 100 -> 40
$ws.Cells.Item(289,19).Value = 12     # S289 ScanDist 10 -> 12
$ws.Cells.Item(289,20).Value = 15     # T289 ChaseDist 12 -> 15
$ws.Cells.Item(289,23).Value = "Holy1" # W289 Element Holy4 -> Holy1
$ws.Cells.Item(289,36).Value = 1.4    # AJ289 ClientSize 1 -> 1.4

# --- Row 290 (existing monster) stat tweaks ---
$ws.Cells.Item(290,5).Value = 75      # E290 HP 100 -> 75
$ws.Cells.Item(290,8).Value = 70      # H290 Vit 100 -> 70
$ws.Cells.Item(290,9).Value = 70      # I290 Dex 100 -> 70
$ws.Cells.Item(290,10).Value = 30     # J290 Agi 100 -> 30
$ws.Cells.Item(290,12).Value = 45     # L290 Attack 50 -> 45
$ws.Cells.Item(290,15).Value = 60     # O290 Def 100 -> 60
$ws.Cells.Item(290,16).Value = 70     # P290 MDef 100 -> 70
$ws.Cells.Item(290,19).Value = 12     # S290 ScanDist 10 -> 12
$ws.Cells.Item(290,20).Value = 15     # T290 ChaseDist 12 -> 15
$ws.Cells.Item(290,23).Value = "Holy1" # W290 Element Holy4 -> Holy1

# --- MonsterAiType fixes (AiAggressive -> AiAngry) ---
$ws.Cells.Item(295,30).Value = "AiAngry"  # AD295
$ws.Cells.Item(299,30).Value = "AiAngry"  # AD299

# --- New rows 320-325 ---
    # Row 320
    $ws.Cells.Item(320,1).Value = 6016
    $ws.Cells.Item(320,2).Value = "SUCCUBUS"
    $ws.Cells.Item(320,3).Value = "Succubus"
    $ws.Cells.Item(320,4).Value = 85
    $ws.Cells.Item(320,5).Value = 100
    $ws.Cells.Item(320,6).Value = 100
    $ws.Cells.Item(320,7).Value = 100
    $ws.Cells.Item(320,8).Value = 100
    $ws.Cells.Item(320,9).Value = 100
    $ws.Cells.Item(320,10).Value = 100
    $ws.Cells.Item(320,11).Value = 100
    $ws.Cells.Item(320,12).Value = 100
    $ws.Cells.Item(320,13).Value = 10
    $ws.Cells.Item(320,14).Value = 2
    $ws.Cells.Item(320,15).Value = 100
    $ws.Cells.Item(320,16).Value = 100
    $ws.Cells.Item(320,17).Value = 100
    $ws.Cells.Item(320,18).Value = 100
    $ws.Cells.Item(320,19).Value = 10
    $ws.Cells.Item(320,20).Value = 12
    $ws.Cells.Item(320,21).Value = "Medium"
    $ws.Cells.Item(320,22).Value = "Demon"
    $ws.Cells.Item(320,23).Value = "Dark3"
    $ws.Cells.Item(320,24).Value = 1306
    $ws.Cells.Item(320,25).Value = 288
    $ws.Cells.Item(320,26).Value = 1056
    $ws.Cells.Item(320,27).Value = 155
    $ws.Cells.Item(320,28).Value = "Normal"
    $ws.Cells.Item(320,29).Value = "Buff,Demon"
    $ws.Cells.Item(320,30).Value = "AiStandardBoss"
    $ws.Cells.Item(320,32).Value = 576
    $ws.Cells.Item(320,33).Value = "succubus.spr"
    $ws.Cells.Item(320,34).Value = 0
    $ws.Cells.Item(320,35).Value = 0.5
    $ws.Cells.Item(320,36).Value = 1
    # Row 321
    $ws.Cells.Item(321,1).Value = 6017
    $ws.Cells.Item(321,2).Value = "INCUBUS"
    $ws.Cells.Item(321,3).Value = "Incubus"
    $ws.Cells.Item(321,4).Value = 75
    $ws.Cells.Item(321,5).Value = 100
    $ws.Cells.Item(321,6).Value = 100
    $ws.Cells.Item(321,7).Value = 100
    $ws.Cells.Item(321,8).Value = 100
    $ws.Cells.Item(321,9).Value = 100
    $ws.Cells.Item(321,10).Value = 100
    $ws.Cells.Item(321,11).Value = 100
    $ws.Cells.Item(321,12).Value = 100
    $ws.Cells.Item(321,13).Value = 10
    $ws.Cells.Item(321,14).Value = 2
    $ws.Cells.Item(321,15).Value = 100
    $ws.Cells.Item(321,16).Value = 100
    $ws.Cells.Item(321,17).Value = 100
    $ws.Cells.Item(321,18).Value = 100
    $ws.Cells.Item(321,19).Value = 10
    $ws.Cells.Item(321,20).Value = 12
    $ws.Cells.Item(321,21).Value = "Medium"
    $ws.Cells.Item(321,22).Value = "Demon"
    $ws.Cells.Item(321,23).Value = "Dark3"
    $ws.Cells.Item(321,24).Value = 850
    $ws.Cells.Item(321,25).Value = 336
    $ws.Cells.Item(321,26).Value = 600
    $ws.Cells.Item(321,27).Value = 165
    $ws.Cells.Item(321,28).Value = "Normal"
    $ws.Cells.Item(321,29).Value = "Buff,Demon"
    $ws.Cells.Item(321,30).Value = "AiStandardBoss"
    $ws.Cells.Item(321,32).Value = 420
    $ws.Cells.Item(321,33).Value = "incubus.spr"
    $ws.Cells.Item(321,34).Value = 0
    $ws.Cells.Item(321,35).Value = 0.5
    $ws.Cells.Item(321,36).Value = 1
    # Row 322
    $ws.Cells.Item(322,1).Value = 6018
    $ws.Cells.Item(322,2).Value = "OBSERVATION"
    $ws.Cells.Item(322,3).Value = "Dame of Sentinel"
    $ws.Cells.Item(322,4).Value = 81
    $ws.Cells.Item(322,5).Value = 100
    $ws.Cells.Item(322,6).Value = 100
    $ws.Cells.Item(322,7).Value = 100
    $ws.Cells.Item(322,8).Value = 100
    $ws.Cells.Item(322,9).Value = 100
    $ws.Cells.Item(322,10).Value = 100
    $ws.Cells.Item(322,11).Value = 100
    $ws.Cells.Item(322,12).Value = 100
    $ws.Cells.Item(322,13).Value = 10
    $ws.Cells.Item(322,14).Value = 2
    $ws.Cells.Item(322,15).Value = 100
    $ws.Cells.Item(322,16).Value = 100
    $ws.Cells.Item(322,17).Value = 100
    $ws.Cells.Item(322,18).Value = 100
    $ws.Cells.Item(322,19).Value = 10
    $ws.Cells.Item(322,20).Value = 12
    $ws.Cells.Item(322,21).Value = "Medium"
    $ws.Cells.Item(322,22).Value = "Angel"
    $ws.Cells.Item(322,23).Value = "Neutral4"
    $ws.Cells.Item(322,24).Value = 432
    $ws.Cells.Item(322,25).Value = 360
    $ws.Cells.Item(322,26).Value = 480
    $ws.Cells.Item(322,27).Value = 100
    $ws.Cells.Item(322,28).Value = "Boss"
    $ws.Cells.Item(322,29).Value = "Strong,Angel"
    $ws.Cells.Item(322,30).Value = "AiStandardBoss"
    $ws.Cells.Item(322,32).Value = 300
    $ws.Cells.Item(322,33).Value = "observation.spr"
    $ws.Cells.Item(322,34).Value = 0
    $ws.Cells.Item(322,35).Value = 0.5
    $ws.Cells.Item(322,36).Value = 1
    # Row 323
    $ws.Cells.Item(323,1).Value = 6019
    $ws.Cells.Item(323,2).Value = "RETRIBUTION"
    $ws.Cells.Item(323,3).Value = "Baroness of Retribution"
    $ws.Cells.Item(323,4).Value = 79
    $ws.Cells.Item(323,5).Value = 100
    $ws.Cells.Item(323,6).Value = 100
    $ws.Cells.Item(323,7).Value = 100
    $ws.Cells.Item(323,8).Value = 100
    $ws.Cells.Item(323,9).Value = 100
    $ws.Cells.Item(323,10).Value = 100
    $ws.Cells.Item(323,11).Value = 100
    $ws.Cells.Item(323,12).Value = 100
    $ws.Cells.Item(323,13).Value = 10
    $ws.Cells.Item(323,14).Value = 2
    $ws.Cells.Item(323,15).Value = 100
    $ws.Cells.Item(323,16).Value = 100
    $ws.Cells.Item(323,17).Value = 100
    $ws.Cells.Item(323,18).Value = 100
    $ws.Cells.Item(323,19).Value = 10
    $ws.Cells.Item(323,20).Value = 12
    $ws.Cells.Item(323,21).Value = "Medium"
    $ws.Cells.Item(323,22).Value = "Angel"
    $ws.Cells.Item(323,23).Value = "Dark3"
    $ws.Cells.Item(323,24).Value = 360
    $ws.Cells.Item(323,25).Value = 360
    $ws.Cells.Item(323,26).Value = 480
    $ws.Cells.Item(323,27).Value = 120
    $ws.Cells.Item(323,28).Value = "Boss"
    $ws.Cells.Item(323,29).Value = "Strong,Angel"
    $ws.Cells.Item(323,30).Value = "AiStandardBoss"
    $ws.Cells.Item(323,32).Value = 240
    $ws.Cells.Item(323,33).Value = "retribution.spr"
    $ws.Cells.Item(323,34).Value = 0
    $ws.Cells.Item(323,35).Value = 0.5
    $ws.Cells.Item(323,36).Value = 1
    # Row 324
    $ws.Cells.Item(324,1).Value = 6020
    $ws.Cells.Item(324,2).Value = "SOLACE"
    $ws.Cells.Item(324,3).Value = "Dame of Sentinel"
    $ws.Cells.Item(324,4).Value = 77
    $ws.Cells.Item(324,5).Value = 100
    $ws.Cells.Item(324,6).Value = 100
    $ws.Cells.Item(324,7).Value = 100
    $ws.Cells.Item(324,8).Value = 100
    $ws.Cells.Item(324,9).Value = 100
    $ws.Cells.Item(324,10).Value = 100
    $ws.Cells.Item(324,11).Value = 100
    $ws.Cells.Item(324,12).Value = 100
    $ws.Cells.Item(324,13).Value = 10
    $ws.Cells.Item(324,14).Value = 2
    $ws.Cells.Item(324,15).Value = 100
    $ws.Cells.Item(324,16).Value = 100
    $ws.Cells.Item(324,17).Value = 100
    $ws.Cells.Item(324,18).Value = 100
    $ws.Cells.Item(324,19).Value = 10
    $ws.Cells.Item(324,20).Value = 12
    $ws.Cells.Item(324,21).Value = "Medium"
    $ws.Cells.Item(324,22).Value = "Angel"
    $ws.Cells.Item(324,23).Value = "Holy3"
    $ws.Cells.Item(324,24).Value = 576
    $ws.Cells.Item(324,25).Value = 360
    $ws.Cells.Item(324,26).Value = 420
    $ws.Cells.Item(324,27).Value = 180
    $ws.Cells.Item(324,28).Value = "Boss"
    $ws.Cells.Item(324,29).Value = "Buff,Angel"
    $ws.Cells.Item(324,30).Value = "AiStandardBoss"
    $ws.Cells.Item(324,32).Value = 384
    $ws.Cells.Item(324,33).Value = "solace.spr"
    $ws.Cells.Item(324,34).Value = 0
    $ws.Cells.Item(324,35).Value = 0.5
    $ws.Cells.Item(324,36).Value = 1
    # Row 325
    $ws.Cells.Item(325,1).Value = 6021
    $ws.Cells.Item(325,2).Value = "SHELTER"
    $ws.Cells.Item(325,3).Value = "Mistress of Shelter"
    $ws.Cells.Item(325,4).Value = 80
    $ws.Cells.Item(325,5).Value = 100
    $ws.Cells.Item(325,6).Value = 100
    $ws.Cells.Item(325,7).Value = 100
    $ws.Cells.Item(325,8).Value = 100
    $ws.Cells.Item(325,9).Value = 100
    $ws.Cells.Item(325,10).Value = 100
    $ws.Cells.Item(325,11).Value = 100
    $ws.Cells.Item(325,12).Value = 100
    $ws.Cells.Item(325,13).Value = 10
    $ws.Cells.Item(325,14).Value = 2
    $ws.Cells.Item(325,15).Value = 100
    $ws.Cells.Item(325,16).Value = 100
    $ws.Cells.Item(325,17).Value = 100
    $ws.Cells.Item(325,18).Value = 100
    $ws.Cells.Item(325,19).Value = 10
    $ws.Cells.Item(325,20).Value = 12
    $ws.Cells.Item(325,21).Value = "Medium"
    $ws.Cells.Item(325,22).Value = "Angel"
    $ws.Cells.Item(325,23).Value = "Holy3"
    $ws.Cells.Item(325,24).Value = 432
    $ws.Cells.Item(325,25).Value = 360
    $ws.Cells.Item(325,26).Value = 420
    $ws.Cells.Item(325,27).Value = 160
    $ws.Cells.Item(325,28).Value = "Boss"
    $ws.Cells.Item(325,29).Value = "Elite,Angel"
    $ws.Cells.Item(325,30).Value = "AiStandardBoss"
    $ws.Cells.Item(325,32).Value = 240
    $ws.Cells.Item(325,33).Value = "shelter.spr"
    $ws.Cells.Item(325,34).Value = 0
    $ws.Cells.Item(325,35).Value = 0.5
    $ws.Cells.Item(325,36).Value = 1

# --- View-state bookkeeping to mirror the diff (best effort; cosmetic) ---
$ws.Activate() | Out-Null
$ws.Range("AD296").Select() | Out-Null

Write-Host "StatDef: updated rows 289, 290, 295, 299; added rows 320-325"
